# Update countries & provincias Spain
#
# The underlying COVID-19 case data was refreshed and the "Pais" sheet
# re-sorted by total cases (column B) descending, as it always is. Most
# rows keep their rank and simply get new figures; two pairs of rows
# swapped rank (Nepal <-> Costa Rica, Malta <-> Sudan del Sur), so their
# country-name cells are updated as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rank swaps: update the country name in column A ---
# Nepal's totals overtook Costa Rica's.
$ws.Range("A55").Value = "Nepal"
$ws.Range("A56").Value = "Costa Rica"
# Malta's totals overtook Sudan del Sur's.
$ws.Range("A145").Value = "Malta"
$ws.Range("A146").Value = "Sudan del Sur"

# --- Refreshed statistics (B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6828698
$ws.Range("C4").Value = 397
$ws.Range("D4").Value = 4119782
$ws.Range("E4").Value = 2507550
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 201366

# Israel (row 27)
$ws.Range("B27").Value = 172322
$ws.Range("C27").Value = 1857
$ws.Range("D27").Value = 125671
$ws.Range("E27").Value = 45488
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 1163

# Oman (row 41)
$ws.Range("B41").Value = 91753
$ws.Range("C41").Value = 557
$ws.Range("D41").Value = 84648
$ws.Range("E41").Value = 6287
$ws.Range("G41").Value = 13
$ws.Range("H41").Value = 818

# Emiratos Arabes Unidos (row 46)
$ws.Range("B46").Value = 82568
$ws.Range("C46").Value = 786
$ws.Range("D46").Value = 72117
$ws.Range("E46").Value = 10049

# Nepal (row 55, formerly Costa Rica's slot)
$ws.Range("B55").Value = 59573
$ws.Range("C55").Value = 1246
$ws.Range("D55").Value = 42949
$ws.Range("E55").Value = 16241
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 383

# Costa Rica (row 56, formerly Nepal's slot) - takes Costa Rica's prior figures
$ws.Range("B56").Value = 59516
$ws.Range("D56").Value = 21752
$ws.Range("E56").Value = 37115
$ws.Range("H56").Value = 649

# Suiza (row 61)
$ws.Range("B61").Value = 48795
$ws.Range("C61").Value = 530
$ws.Range("E61").Value = 6854
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 2041

# Estado de Palestina (row 71)
$ws.Range("B71").Value = 33843
$ws.Range("C71").Value = 837
$ws.Range("D71").Value = 23060
$ws.Range("E71").Value = 10539
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 244

# Senegal (row 88)
$ws.Range("B88").Value = 14618
$ws.Range("C88").Value = 50
$ws.Range("D88").Value = 10923
$ws.Range("E88").Value = 3395
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 300

# Consejo Danes para los Refugiados (row 95)
$ws.Range("B95").Value = 10442
$ws.Range("C95").Value = 28
$ws.Range("D95").Value = 9840
$ws.Range("E95").Value = 335

# Malta (row 145, formerly Sudan del Sur's slot)
$ws.Range("B145").Value = 2595
$ws.Range("C145").Value = 35
$ws.Range("D145").Value = 1978
$ws.Range("E145").Value = 601
$ws.Range("H145").Value = 16

# Sudan del Sur (row 146, formerly Malta's slot) - takes Sudan del Sur's prior figures
$ws.Range("B146").Value = 2594
$ws.Range("D146").Value = 1290
$ws.Range("E146").Value = 1255
$ws.Range("H146").Value = 49

# Vietnam (row 168)
$ws.Range("B168").Value = 1066
$ws.Range("C168").Value = 3
$ws.Range("D168").Value = 940
$ws.Range("E168").Value = 91

# Gibraltar (row 182)
$ws.Range("B182").Value = 343
$ws.Range("C182").Value = 3
$ws.Range("D182").Value = 313
$ws.Range("E182").Value = 30

# --- "Last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 13:41"
